$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 12:24"

# Update country data rows (country name swaps + refreshed case numbers)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 2208486
$ws.Cells.Item(4, 3).Value = 86
$ws.Cells.Item(4, 4).Value = 903136
$ws.Cells.Item(4, 5).Value = 1186217
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 119133

$ws.Cells.Item(25, 1).Value = "Belgica"
$ws.Cells.Item(25, 2).Value = 60244
$ws.Cells.Item(25, 3).Value = 89
$ws.Cells.Item(25, 4).Value = 16684
$ws.Cells.Item(25, 5).Value = 33885
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 12
$ws.Cells.Item(25, 8).Value = 9675

$ws.Cells.Item(26, 1).Value = "Bielorrusia"
$ws.Cells.Item(26, 2).Value = 56032
$ws.Cells.Item(26, 3).Value = 663
$ws.Cells.Item(26, 4).Value = 32735
$ws.Cells.Item(26, 5).Value = 22973
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 6
$ws.Cells.Item(26, 8).Value = 324

$ws.Cells.Item(33, 1).Value = "Indonesia"
$ws.Cells.Item(33, 2).Value = 41431
$ws.Cells.Item(33, 3).Value = 1031
$ws.Cells.Item(33, 4).Value = 16243
$ws.Cells.Item(33, 5).Value = 22912
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 45
$ws.Cells.Item(33, 8).Value = 2276

$ws.Cells.Item(34, 1).Value = "Singapur"
$ws.Cells.Item(34, 2).Value = 41216
$ws.Cells.Item(34, 3).Value = 247
$ws.Cells.Item(34, 4).Value = 31163
$ws.Cells.Item(34, 5).Value = 10027
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 26

$ws.Cells.Item(46, 1).Value = "Rumania"
$ws.Cells.Item(46, 2).Value = 22760
$ws.Cells.Item(46, 3).Value = 345
$ws.Cells.Item(46, 4).Value = 16117
$ws.Cells.Item(46, 5).Value = 5192
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 14
$ws.Cells.Item(46, 8).Value = 1451

$ws.Cells.Item(47, 1).Value = "Irak"
$ws.Cells.Item(47, 2).Value = 22700
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 9862
$ws.Cells.Item(47, 5).Value = 12126
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 712

$ws.Cells.Item(68, 1).Value = "Marruecos"
$ws.Cells.Item(68, 2).Value = 8985
$ws.Cells.Item(68, 3).Value = 54
$ws.Cells.Item(68, 4).Value = 7960
$ws.Cells.Item(68, 5).Value = 813
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 212

$ws.Cells.Item(70, 1).Value = "Malasia"
$ws.Cells.Item(70, 2).Value = 8515
$ws.Cells.Item(70, 3).Value = 10
$ws.Cells.Item(70, 4).Value = 7873
$ws.Cells.Item(70, 5).Value = 521
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 121

$ws.Cells.Item(73, 1).Value = "Finlandia"
$ws.Cells.Item(73, 2).Value = 7117
$ws.Cells.Item(73, 3).Value = 5
$ws.Cells.Item(73, 4).Value = 6200
$ws.Cells.Item(73, 5).Value = 591
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 326

$ws.Cells.Item(79, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(79, 2).Value = 5100
$ws.Cells.Item(79, 3).Value = 126
$ws.Cells.Item(79, 4).Value = 640
$ws.Cells.Item(79, 5).Value = 4345
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 3
$ws.Cells.Item(79, 8).Value = 115

$ws.Cells.Item(104, 1).Value = "Sri Lanka"
$ws.Cells.Item(104, 2).Value = 1915
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 1397
$ws.Cells.Item(104, 5).Value = 507
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 11

$ws.Cells.Item(111, 1).Value = "Albania"
$ws.Cells.Item(111, 2).Value = 1722
$ws.Cells.Item(111, 3).Value = 50
$ws.Cells.Item(111, 4).Value = 1077
$ws.Cells.Item(111, 5).Value = 607
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 38

$ws.Cells.Item(115, 1).Value = "Eslovenia"
$ws.Cells.Item(115, 2).Value = 1503
$ws.Cells.Item(115, 3).Value = 4
$ws.Cells.Item(115, 4).Value = 1359
$ws.Cells.Item(115, 5).Value = 35
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 109

$ws.Cells.Item(137, 1).Value = "Uganda"
$ws.Cells.Item(137, 2).Value = 732
$ws.Cells.Item(137, 3).Value = 8
$ws.Cells.Item(137, 4).Value = 420
$ws.Cells.Item(137, 5).Value = 312
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 0

$ws.Cells.Item(206, 1).Value = "Groenlandia"
$ws.Cells.Item(206, 2).Value = 13
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 13
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

$ws.Cells.Item(207, 1).Value = "Islas Malvinas"
$ws.Cells.Item(207, 2).Value = 13
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 13
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 1

$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 2).Value = 8
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1
